# v2 rollout - "clean and slim version"
# - Rename "Class Quizes" sheet to "assessment"
# - Add a defined name "edem" pointing at assessment!$M$29
# - Slightly widen a couple of label columns
# - Drop the now-unused trailing row (and its "\" label) from the assessment sheet

$wb = $excel.ActiveWorkbook

$wsPresentation = $wb.Worksheets.Item("Presentation Groups")
$wsAssessment   = $wb.Worksheets.Item("Class Quizes")

# Rename "Class Quizes" -> "assessment"
$wsAssessment.Name = "assessment"

# New named range "edem" referring to the running total cell on the assessment sheet
# (single-quoted so PowerShell doesn't try to expand $M as a variable)
$wb.Names.Add("edem", '=assessment!$M$29')

# Widen "Presentation Groups" column B slightly
$wsPresentation.Columns.Item(2).ColumnWidth = 43.666666666666664

# Widen the assessment sheet's data columns (B through L) slightly
$wsAssessment.Columns.Item(2).ColumnWidth  = 38.833333333333336
$wsAssessment.Columns.Item(3).ColumnWidth  = 23.666666666666668
$wsAssessment.Columns.Item(4).ColumnWidth  = 24.166666666666668
$wsAssessment.Columns.Item(5).ColumnWidth  = 25.5
$wsAssessment.Columns.Item(6).ColumnWidth  = 25.166666666666668
$wsAssessment.Columns.Item(7).ColumnWidth  = 19.833333333333332
$wsAssessment.Columns.Item(8).ColumnWidth  = 21.833333333333332
$wsAssessment.Columns.Item(9).ColumnWidth  = 23.5
$wsAssessment.Columns.Item(10).ColumnWidth = 24.833333333333332
$wsAssessment.Columns.Item(11).ColumnWidth = 22.5
$wsAssessment.Columns.Item(12).ColumnWidth = 28.5

# Move the view/selection to where the editor left off, then drop the stray
# trailing row 54 (which only held a lone "\" marker in column L)
$wsAssessment.Activate()
$wsAssessment.Range("J52").Select()
$wsAssessment.Rows.Item(54).Delete()
